$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Samples" query text for B3 (SamplesTab query) ---
# same as before, but with a double space before the predicate list
# ("IN  [" instead of "IN [")
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN  ["Taxane only"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
 order By samp.sample_id ASC LIMIT 100
'@

# --- New shared "StatQuery" text used by both the Samples and Files tabs (C3 & C4) ---
$statQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen in  ["Taxane only"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

# Row 3 = SamplesTab : update the stat query (C3) and the query (B3)
# (stat query is assigned first so the new shared-string table keeps the
# same relative ordering as the source workbook)
$ws.Range("C3").Value = $statQuery
$ws.Range("B3").Value = $samplesQuery

# Row 4 = FilesTab : query (B4) is unchanged, only the stat query (C4) is updated
$ws.Range("C4").Value = $statQuery

# New row 5: a single formatted (wrap-text) empty cell in column C,
# matching the style used by the query cells
$ws.Range("C5").Value = ""
$ws.Range("C5").WrapText = $true

# Row heights were adjusted slightly when the workbook was last saved
$ws.Rows.Item(2).RowHeight = 345
$ws.Rows.Item(3).RowHeight = 375
$ws.Rows.Item(4).RowHeight = 409.5

# Update view/selection to match the saved state (scrolled up one row,
# with B3 selected instead of B4)
[void]$ws.Range("B3").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
